$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.371.77"
$ws.Cells.Item(2, 5).Value = "  -2.92%  "

$ws.Cells.Item(3, 4).Value = "2.244.67"
$ws.Cells.Item(3, 5).Value = "  -3.82%  "

$ws.Cells.Item(4, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(5, 4)
$c.Value = "'233.90"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.86%  "

$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.632"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -4.26%  "

$ws.Cells.Item(7, 5).Value = "  -2.77%  "

$ws.Cells.Item(8, 5).Value = "  +0.01%  "

$ws.Cells.Item(9, 5).Value = "  -3.32%  "

$ws.Cells.Item(10, 5).Value = "  +1.48%  "

$c = $ws.Cells.Item(11, 4)
$c.Value = "'58.79"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.38%  "

$c = $ws.Cells.Item(12, 4)
$c.Value = "'36.55"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +13.37%  "

$c = $ws.Cells.Item(13, 4)
$c.Value = "'0.106"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.50%  "

$ws.Cells.Item(14, 5).Value = "  -4.89%  "

$ws.Cells.Item(15, 4).Value = "2.579.90"
$ws.Cells.Item(15, 5).Value = "  -3.75%  "

$c = $ws.Cells.Item(16, 4)
$c.Value = "'15.10"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -5.30%  "

$c = $ws.Cells.Item(17, 4)
$c.Value = "'0.857"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.54%  "

$ws.Cells.Item(18, 4).Value = "2.248.20"
$ws.Cells.Item(18, 5).Value = "  -3.88%  "

$ws.Cells.Item(19, 4).Value = "42.267.68"
$ws.Cells.Item(19, 5).Value = "  -3.05%  "

$ws.Cells.Item(20, 5).Value = "  -2.33%  "

$c = $ws.Cells.Item(21, 4)
$c.Value = "'6.25"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -4.60%  "

$c = $ws.Cells.Item(22, 4)
$c.Value = "'73.45"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -5.61%  "

$c = $ws.Cells.Item(23, 4)
$c.Value = "'236.31"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -5.41%  "

$c = $ws.Cells.Item(24, 4)
$c.Value = "'1.98"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +4.45%  "

$ws.Cells.Item(25, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(26, 4)
$c.Value = "'3.68"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.63%  "

$ws.Cells.Item(27, 5).Value = "  -3.37%  "

$c = $ws.Cells.Item(28, 4)
$c.Value = "'10.00"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -2.35%  "

$ws.Cells.Item(29, 5).Value = "  -1.99%  "

$c = $ws.Cells.Item(30, 4)
$c.Value = "'170.16"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -3.07%  "

$c = $ws.Cells.Item(31, 4)
$c.Value = "'20.56"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -6.47%  "

$ws.Cells.Item(32, 5).Value = "  -3.38%  "

$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.126"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.17%  "

$ws.Cells.Item(34, 5).Value = "  +0.51%  "

$ws.Cells.Item(35, 5).Value = "  +0.65%  "

$c = $ws.Cells.Item(36, 4)
$c.Value = "'4.71"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -6.15%  "

$c = $ws.Cells.Item(37, 4)
$c.Value = "'3.67"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.11%  "

$ws.Cells.Item(38, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(38, 4)
$c.Value = "'22.03"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +18.23%  "

$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.0281"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +4.75%  "

$ws.Cells.Item(40, 5).Value = "  -2.96%  "

$c = $ws.Cells.Item(41, 4)
$c.Value = "'5.99"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -5.66%  "

$c = $ws.Cells.Item(42, 4)
$c.Value = "'65.40"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.80%  "

$c = $ws.Cells.Item(43, 4)
$c.Value = "'9.22"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.92%  "

$c = $ws.Cells.Item(44, 4)
$c.Value = "'4.94"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -11.91%  "

$ws.Cells.Item(45, 5).Value = "  -2.10%  "

$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.191"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.30%  "

$ws.Cells.Item(47, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(48, 4)
$c.Value = "'4.51"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +11.15%  "

$c = $ws.Cells.Item(49, 4)
$c.Value = "'10.21"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +10.27%  "

$ws.Cells.Item(50, 5).Value = "  -2.60%  "

$c = $ws.Cells.Item(51, 4)
$c.Value = "'2.34"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.11%  "
